$d = $word.ActiveDocument

# The document currently ends with the "Week 8" paragraph, immediately
# followed by the section properties. Append a brand new paragraph after
# it for the "Week 9" log entry.

$lastPara = $d.Paragraphs.Last
$tailRange = $lastPara.Range
$tailRange.InsertParagraphAfter() | Out-Null

# The newly created paragraph is now the last paragraph in the document;
# it inherits the Tahoma run formatting from the preceding "Week 8"
# paragraph automatically.
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
# Exclude the trailing paragraph mark so inserted text lands inside the
# paragraph (before the pilcrow), matching how Word itself behaves.
$newRange.End = $newRange.End - 1

# Matches the author's several back-to-back edits/pastes that together
# make up the Week 9 entry (each becomes its own run in the OOXML).
$segments = @(
  "Week 9: Done doing a trial run trying to connect the BLE to the phone using an existing ",
  "ESP32 BLE server ",
  "library codes and it works. But to connect it, you have to use an app that can connect to the ESP32 because using the local settings in the iPhone does not connect to the ESP32. After doing that I went straight into watching tutorials on JavaScript to prepare myself for node.js for the app making.",
  " I’m going to finish watching",
  " JavaScript",
  " tutorials until the end of the week 9."
)

$fullText = [string]::Join("", $segments)
$newRange.Text = $fullText

# Touch (and immediately revert) a throwaway formatting property on each
# segment's own sub-range. This keeps the final, visible formatting
# untouched (still just the inherited Tahoma rFonts) while still forcing
# each segment to remain its own distinct run in the saved document,
# rather than being silently coalesced with its neighbours.
$pos = $newRange.Start
foreach ($seg in $segments) {
  $segRange = $d.Range($pos, $pos + $seg.Length)
  $segRange.Bold = 1
  $segRange.Bold = 0
  $pos = $pos + $seg.Length
}
